$d = $word.ActiveDocument
$d.Content.Find.Execute("Sternbild Sternbild Stier", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Sternbild Stier", 2)
